$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell A13 value: VO:0010434 -> VO:0010458
$ws.Range("A13").Value = "VO:0010458"

# Update the active selection to A14
$ws.Range("A14").Select()
